{"js": "// Update the worksheet date and the 25 two-digit multiplication problems\n// (aa3dc9e): each old \"A\u00d7B=\" expression is replaced by its new value, and\n// the heading date line is bumped from 2024-01-07 Sunday to 2024-01-08 Monday.\n\nconst replacements = [\n    [\"2024-01-07 Sunday\", \"2024-01-08 Monday\"],\n    [\"58\u00d740=\", \"69\u00d783=\"],\n    [\"71\u00d720=\", \"23\u00d799=\"],\n    [\"78\u00d742=\", \"53\u00d785=\"],\n    [\"40\u00d747=\", \"53\u00d794=\"],\n    [\"54\u00d718=\", \"90\u00d758=\"],\n    [\"33\u00d767=\", \"68\u00d767=\"],\n    [\"70\u00d712=\", \"87\u00d717=\"],\n    [\"28\u00d717=\", \"84\u00d799=\"],\n    [\"79\u00d766=\", \"95\u00d714=\"],\n    [\"64\u00d759=\", \"38\u00d773=\"],\n    [\"97\u00d756=\", \"32\u00d774=\"],\n    [\"88\u00d725=\", \"61\u00d717=\"],\n    [\"30\u00d715=\", \"87\u00d782=\"],\n    [\"70\u00d730=\", \"32\u00d750=\"],\n    [\"20\u00d730=\", \"99\u00d779=\"],\n    [\"71\u00d749=\", \"56\u00d748=\"],\n    [\"16\u00d781=\", \"59\u00d735=\"],\n    [\"54\u00d738=\", \"60\u00d717=\"],\n    [\"84\u00d782=\", \"37\u00d745=\"],\n    [\"96\u00d723=\", \"59\u00d717=\"],\n    [\"16\u00d750=\", \"22\u00d725=\"],\n    [\"58\u00d714=\", \"31\u00d761=\"],\n    [\"62\u00d747=\", \"50\u00d728=\"],\n    [\"26\u00d790=\", \"98\u00d769=\"],\n    [\"34\u00d726=\", \"85\u00d750=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit multiplication problems\n# (aa3dc9e): each old \"A\u00d7B=\" expression is replaced by its new value, and\n# the heading date line is bumped from 2024-01-07 Sunday to 2024-01-08 Monday.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-07 Sunday\", \"2024-01-08 Monday\"),\n    @(\"58\u00d740=\", \"69\u00d783=\"),\n    @(\"71\u00d720=\", \"23\u00d799=\"),\n    @(\"78\u00d742=\", \"53\u00d785=\"),\n    @(\"40\u00d747=\", \"53\u00d794=\"),\n    @(\"54\u00d718=\", \"90\u00d758=\"),\n    @(\"33\u00d767=\", \"68\u00d767=\"),\n    @(\"70\u00d712=\", \"87\u00d717=\"),\n    @(\"28\u00d717=\", \"84\u00d799=\"),\n    @(\"79\u00d766=\", \"95\u00d714=\"),\n    @(\"64\u00d759=\", \"38\u00d773=\"),\n    @(\"97\u00d756=\", \"32\u00d774=\"),\n    @(\"88\u00d725=\", \"61\u00d717=\"),\n    @(\"30\u00d715=\", \"87\u00d782=\"),\n    @(\"70\u00d730=\", \"32\u00d750=\"),\n    @(\"20\u00d730=\", \"99\u00d779=\"),\n    @(\"71\u00d749=\", \"56\u00d748=\"),\n    @(\"16\u00d781=\", \"59\u00d735=\"),\n    @(\"54\u00d738=\", \"60\u00d717=\"),\n    @(\"84\u00d782=\", \"37\u00d745=\"),\n    @(\"96\u00d723=\", \"59\u00d717=\"),\n    @(\"16\u00d750=\", \"22\u00d725=\"),\n    @(\"58\u00d714=\", \"31\u00d761=\"),\n    @(\"62\u00d747=\", \"50\u00d728=\"),\n    @(\"26\u00d790=\", \"98\u00d769=\"),\n    @(\"34\u00d726=\", \"85\u00d750=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
